$d = $word.ActiveDocument

# Update the date heading (unique text, safe to use Find & Replace)
$d.Content.Find.Execute("2025-07-01 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-07-02 Wednesday", 2) | Out-Null

# Update the division problems in the table. Cell coordinates are used
# (rather than a global text Find & Replace) because several of the
# "before" and "after" values collide with each other across cells,
# and per-cell assignment avoids any double-substitution issues.
$t = $d.Tables.Item(1)

$newValues = @{
    "1,1" = "31÷8=";  "1,2" = "37÷8=";  "1,3" = "99÷7=";  "1,4" = "93÷8=";  "1,5" = "36÷8=";
    "5,1" = "49÷9=";  "5,2" = "68÷2=";  "5,3" = "86÷8=";  "5,4" = "31÷7=";  "5,5" = "83÷7=";
    "9,1" = "24÷8=";  "9,2" = "61÷7=";  "9,3" = "98÷4=";  "9,4" = "80÷5=";  "9,5" = "15÷5=";
    "13,1" = "13÷9="; "13,2" = "76÷3="; "13,3" = "89÷2="; "13,4" = "24÷3="; "13,5" = "76÷5=";
    "17,1" = "77÷5="; "17,2" = "47÷2="; "17,3" = "77÷8="; "17,4" = "19÷9="; "17,5" = "34÷9=";
}

foreach ($row in @(1, 5, 9, 13, 17)) {
    for ($col = 1; $col -le 5; $col++) {
        $key = "$row,$col"
        $cell = $t.Cell($row, $col)
        $cell.Range.Text = $newValues[$key]
    }
}
